$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the SQL queries embedded in column B/C: the join conditions were
# widened from the bare "id" columns to the fully-qualified
# "<table>_id" / "<table>.<table>_id" forms.
$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $cells) {
    $text = $ws.Range($addr).Value2
    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    $ws.Range($addr).Value2 = $text
}

# Column C is wider now that the queries are longer (and no longer auto-fit).
$ws.Columns.Item(3).ColumnWidth = 68

# Selection/scroll moved back to the top of the sheet.
$ws.Range("B2").Select()
